$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.870.44"
$ws.Range("E2").Value = "  +0.94%  "

$ws.Range("D3").Value = "2.635.10"
$ws.Range("E3").Value = "  +1.76%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.37"
$ws.Range("E5").Value = "  +3.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.19"
$ws.Range("E6").Value = "  +1.16%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.65"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("E10").Value = "  +5.14%  "

$ws.Range("E11").Value = "  +1.51%  "

$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").Value = "3.096.37"
$ws.Range("E13").Value = "  +1.71%  "

$ws.Range("D14").Value = "60.875.43"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.95"
$ws.Range("E15").Value = "  +2.19%  "

$ws.Range("E16").Value = "  +3.03%  "

$ws.Range("D17").Value = "2.647.13"
$ws.Range("E17").Value = "  +2.18%  "

$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "353.68"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.62"
$ws.Range("E20").Value = "  +1.24%  "

$ws.Range("E21").Value = "  +1.85%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.57"
$ws.Range("E23").Value = "  +2.13%  "

$ws.Range("E24").Value = "  +2.52%  "

$ws.Range("E25").Value = "  +1.29%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").Value = "0.0₃0866"
$ws.Range("E27").Value = "  +3.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  +0.95%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.17"
$ws.Range("E30").Value = "  +7.98%  "

$ws.Range("E31").Value = "  +0.61%  "

$ws.Range("E32").Value = "  +3.95%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.00"
$ws.Range("E33").Value = "  -1.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.14"
$ws.Range("E34").Value = "  +4.27%  "

$ws.Range("E35").Value = "  +2.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.928"
$ws.Range("E36").Value = "  +10.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.893"
$ws.Range("E37").Value = "  +3.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("E38").Value = "  +1.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.55"
$ws.Range("E41").Value = "  +1.26%  "

$ws.Range("E42").Value = "  +3.91%  "

$ws.Range("E43").Value = "  +1.48%  "

$ws.Range("E44").Value = "  +1.88%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.75"
$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.93"
$ws.Range("E47").Value = "  +3.04%  "

$ws.Range("E48").Value = "  +2.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.33"
$ws.Range("E49").Value = "  +8.55%  "

$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("D51").Value = "1.982.38"
$ws.Range("E51").Value = "  -0.25%  "

# Row 39 and 40: Bittensor/Filecoin swap places with updated values
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "306.43"
$ws.Range("E39").Value = "  +4.24%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.82"
$ws.Range("E40").Value = "  +1.73%  "
